$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Add the new "Colophon: text" paragraph style (styleId "Colophontext"),
#    based on Normal, quick style.
# ---------------------------------------------------------------------
$colophonStyle = $d.Styles.Add("Colophon: text", 1)
$colophonStyle.BaseStyle = $d.Styles("Normal")
$colophonStyle.QuickStyle = $true

# ---------------------------------------------------------------------
# 2. Apply the new style to the colophon paragraphs.
#    Paragraphs 29-32 keep their existing centre alignment; 33-37 just
#    get the style (33/34/36 are the blank spacer paragraphs, 35/37 keep
#    their existing - non centred - alignment).
# ---------------------------------------------------------------------
foreach ($i in 29, 30, 31, 32) {
    $p = $d.Paragraphs($i)
    $p.Style = "Colophon: text"
    $p.Alignment = 1
}
foreach ($i in 33, 34, 35, 36, 37) {
    $p = $d.Paragraphs($i)
    $p.Style = "Colophon: text"
}

# ---------------------------------------------------------------------
# 3. Split the "© " run into separate "©" and " " runs.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("© ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightRng = $d.Range($rng.Start, $rng.Start + 1)
$copyrightRng.Bold = 1
$copyrightRng.Bold = 0

# ---------------------------------------------------------------------
# 4. Catalogue number / ISBN / ISSN paragraph updates.
# ---------------------------------------------------------------------

# Split "N° de cat. " into "N° de cat." + " "
$rng = $d.Content
$rng.Find.Execute("N° de cat. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$catRng = $d.Range($rng.Start, $rng.End - 1)
$catRng.Bold = 1
$catRng.Bold = 0

# Update the catalogue number bookmark text
$rng = $d.Content
$rng.Find.Execute("Fs97-6/0F-PDF", $true, $false, $false, $false, $false, $true, 1, $false, "Fs97-6/3718E-PDF", 2)

# Merge the "     " + "ISBN " runs into a single run
$rng = $d.Content
$rng.Find.Execute("     ISBN ", $true, $false, $false, $false, $false, $true, 1, $false, "     ISBN ", 2)

# Update the ISBN bookmark text
$rng = $d.Content
$rng.Find.Execute("978-0-660-73817-8", $true, $false, $false, $false, $false, $true, 1, $false, "978-0-660-78637-7", 2)

# Merge the "     " + "ISSN 1488-545X" runs into a single run and update the number
$rng = $d.Content
$rng.Find.Execute("     ISSN 1488-545X", $true, $false, $false, $false, $false, $true, 1, $false, "     ISSN 1488-5379", 2)

# ---------------------------------------------------------------------
# 5. "suit" + " " -> single run "suit ", keeping ":" as its own run.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("suit ", $true, $false, $false, $false, $false, $true, 1, $false, "suit ", 2)

$rng = $d.Content
$rng.Find.Execute("suit :", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$colonRng = $d.Range($rng.End - 1, $rng.End)
$colonRng.Bold = 1
$colonRng.Bold = 0

Write-Output "done"
